$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '57.911.57'
$ws.Range("E2").Value = '  +2.76%  '

$ws.Range("D3").Value = '2.359.39'
$ws.Range("E3").Value = '  +1.57%  '

$ws.Range("D4").Value = "'0.998"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.23%  '

$ws.Range("D5").Value = "'544.03"
$ws.Range("D5").Style = "Normal"

$ws.Range("D6").Value = "'134.90"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.57%  '

$ws.Range("E7").Value = '  -0.07%  '

$ws.Range("D8").Value = "'0.536"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.75%  '

$ws.Range("D9").Value = '2.357.53'
$ws.Range("E9").Value = '  +1.40%  '

$ws.Range("E10").Value = '  +1.78%  '

$ws.Range("E11").Value = '  +1.32%  '

$ws.Range("E12").Value = '  +2.79%  '

$ws.Range("D13").Value = "'0.358"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +6.48%  '

$ws.Range("D14").Value = '2.777.66'
$ws.Range("E14").Value = '  +1.44%  '

$ws.Range("D15").Value = "'23.60"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.65%  '

$ws.Range("D16").Value = '57.970.20'
$ws.Range("E16").Value = '  +2.90%  '

$ws.Range("E17").Value = '  +1.40%  '

$ws.Range("D18").Value = '2.356.29'
$ws.Range("E18").Value = '  +1.37%  '

$ws.Range("D19").Value = "'10.59"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.72%  '

$ws.Range("D20").Value = "'335.36"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +3.09%  '

$ws.Range("D21").Value = "'4.21"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.37%  '

$ws.Range("D22").Value = "'6.73"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.16%  '

$ws.Range("D23").Value = "'0.999"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.04%  '

$ws.Range("D24").Value = "'61.81"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.40%  '

$ws.Range("E25").Value = '  +4.75%  '

$ws.Range("B26").Value = 'Binance-PegBSC-USD'
$ws.Range("C26").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D26").Value = "'1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.01%  '

$ws.Range("B27").Value = 'InternetComputer(DFINITY)'
$ws.Range("C27").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D27").Value = "'8.43"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -3.25%  '

$ws.Range("D28").Value = "'1.43"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +9.73%  '

$ws.Range("E29").Value = '  +5.58%  '

$ws.Range("D30").Value = "'169.69"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.44%  '

$ws.Range("E31").Value = '  +3.07%  '

$ws.Range("E32").Value = '  +0.98%  '

$ws.Range("B33").Value = 'EthereumClassic'
$ws.Range("C33").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D33").Value = "'18.54"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.29%  '

$ws.Range("B34").Value = 'SuiNetwork'
$ws.Range("C34").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D34").Value = "'1.03"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +16.03%  '

$ws.Range("D35").Value = "'0.999"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.00%  '

$ws.Range("D36").Value = "'0.996"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.18%  '

$ws.Range("D37").Value = "'4.19"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +7.47%  '

$ws.Range("E38").Value = '  +0.82%  '

$ws.Range("D39").Value = "'1.63"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +5.62%  '

$ws.Range("D40").Value = "'39.33"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.42%  '

$ws.Range("D41").Value = "'149.92"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.69%  '

$ws.Range("E42").Value = '  +1.95%  '

$ws.Range("E43").Value = '  +2.29%  '

$ws.Range("D44").Value = "'287.78"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +4.15%  '

$ws.Range("D45").Value = "'19.37"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +7.65%  '

$ws.Range("D46").Value = "'0.0929"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.33%  '

$ws.Range("D47").Value = "'0.0507"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.83%  '

$ws.Range("D48").Value = "'0.562"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.08%  '

$ws.Range("E49").Value = '  +2.50%  '

$ws.Range("D50").Value = "'17.63"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +3.63%  '

$ws.Range("E51").Value = '  +2.31%  '
